$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Keegan Murray"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Sacramento Kings"

$ws.Range("A10").Value = "Jarrett Allen"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Cleveland Cavaliers"

$ws.Range("A11").Value = "Daniel Gafford"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Dallas Mavericks"

$ws.Range("A13").Value = "Kevin Durant"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Phoenix Suns"

$ws.Range("A14").Value = "Karl-Anthony Towns"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "New York Knicks"

$ws.Range("A15").Value = "Jalen Johnson"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Atlanta Hawks"

$ws.Range("A18").Value = "Mark Williams"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Charlotte Hornets"
